$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Keahlian ID" -> "Bidang Keahlian"
$ws.Range("A1").Value = "Bidang Keahlian"

# The "Keahlian" column now stores the skill names for the row instead of an ID
$ws.Range("A2").Value = "Frontend Developer, Backend Developer"

# Widen column A to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 44.703

$wb.Save()
